# Update the "取得日時" (acquired datetime) column A for rows 2-10
# from "2025-12-30 01:25:08" to "2025-12-30 01:58:06"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newValue = "2025-12-30 01:58:06"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newValue
}
